# Apply updated cryptocurrency price/label values to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    # Force text storage (matches the source data's original inline-string
    # type) so Excel's automatic type detection doesn't coerce the
    # numeric-looking strings into actual numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Column D (Price) updates
Set-TextValue "D2"  "243.82"
Set-TextValue "D3"  "23.97"
Set-TextValue "D4"  "5.105"
Set-TextValue "D5"  "0.05762"
Set-TextValue "D7"  "3.149"
Set-TextValue "D8"  "0.8114"
Set-TextValue "D9"  "0.8426"
Set-TextValue "D10" "0.1349"
Set-TextValue "D11" "0.06957"
Set-TextValue "D12" "0.03119"
Set-TextValue "D13" "0.02843"
Set-TextValue "D14" "0.09373"
Set-TextValue "D15" "3.759"
Set-TextValue "D16" "0.001511"
Set-TextValue "D17" "0.04684"
Set-TextValue "D18" "0.0005981"
Set-TextValue "D19" "0.006143"
Set-TextValue "D20" "0.001240"
Set-TextValue "D22" "0.00008716"
Set-TextValue "D40" "0.03622"
Set-TextValue "D41" "0.006312"
Set-TextValue "D44" "0.007386"
Set-TextValue "D45" "0.00005292"
Set-TextValue "D47" "0.3006"
Set-TextValue "D48" "0.002272"
Set-TextValue "D49" "0.00002104"
Set-TextValue "D50" "0.0002004"

# Column E (Volume(1h)) updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
